# Updates Gungnir_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with refreshed market-board price/profit data from the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 5579157.5
$ws.Range("I40").Value = 10417907
$ws.Range("J40").Value = 1431657.4
$ws.Range("K40").Value = 10417907
$ws.Range("L40").Value = 1431657.4
$ws.Range("M40").Value = -10417732
$ws.Range("N40").Value = -1432007.4

# Row 137
$ws.Range("H137").Value = 2354.2239
$ws.Range("I137").Value = 2132.42
$ws.Range("J137").Value = 3006.5881
$ws.Range("K137").Value = 6397.26
$ws.Range("L137").Value = 9019.764299999999
$ws.Range("M137").Value = -3847.26
$ws.Range("N137").Value = -14119.7643

# Row 138
$ws.Range("H138").Value = 3978.7144
$ws.Range("I138").Value = 1496.081
$ws.Range("J138").Value = 6762.273
$ws.Range("K138").Value = 4488.242999999999
$ws.Range("L138").Value = 20286.819
$ws.Range("M138").Value = 651.7570000000005
$ws.Range("N138").Value = -30566.819

# Row 141
$ws.Range("H141").Value = 2564.25
$ws.Range("I141").Value = 1317.1538
$ws.Range("J141").Value = 7968.3335
$ws.Range("K141").Value = 3951.4614
$ws.Range("L141").Value = 23905.0005
$ws.Range("M141").Value = 1228.5386
$ws.Range("N141").Value = -34265.00049999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 15000
$ws.Range("J6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("N6").Value = -15346

# Row 32
$ws.Range("H32").Value = 5133389
$ws.Range("I32").Value = 4530.9644
$ws.Range("J32").Value = 37046284
$ws.Range("K32").Value = 4530.9644
$ws.Range("L32").Value = 37046284
$ws.Range("M32").Value = -4243.9644
$ws.Range("N32").Value = -37046858

# Row 37
$ws.Range("H37").Value = 12000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 12000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 12000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -12546

# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# Row 74
$ws.Range("H74").Value = 1512.3478
$ws.Range("I74").Value = 1016.7143
$ws.Range("J74").Value = 2283.3333
$ws.Range("K74").Value = 1016.7143
$ws.Range("L74").Value = 2283.3333
$ws.Range("M74").Value = -142.7143
$ws.Range("N74").Value = -4031.3333

# Row 77
$ws.Range("H77").Value = 1512.3478
$ws.Range("I77").Value = 1016.7143
$ws.Range("J77").Value = 2283.3333
$ws.Range("K77").Value = 5083.5715
$ws.Range("L77").Value = 11416.6665
$ws.Range("M77").Value = -715.5715
$ws.Range("N77").Value = -20152.6665

# Row 109
$ws.Range("H109").Value = 54000
$ws.Range("J109").Value = 54000
$ws.Range("L109").Value = 54000
$ws.Range("N109").Value = -56774

# Row 120
$ws.Range("H120").Value = 41399.5
$ws.Range("J120").Value = 41399.5
$ws.Range("L120").Value = 41399.5
$ws.Range("N120").Value = -51075.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1301.1351
$ws.Range("I86").Value = 1376.3928
$ws.Range("J86").Value = 1067
$ws.Range("K86").Value = 1376.3928
$ws.Range("L86").Value = 1067
$ws.Range("M86").Value = -253.3928000000001
$ws.Range("N86").Value = -3313

# Row 89
$ws.Range("H89").Value = 1301.1351
$ws.Range("I89").Value = 1376.3928
$ws.Range("J89").Value = 1067
$ws.Range("K89").Value = 6881.964
$ws.Range("L89").Value = 5335
$ws.Range("M89").Value = -1265.964
$ws.Range("N89").Value = -16567

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 51
$ws.Range("H51").Value = 33057.5
$ws.Range("I51").Value = 2090
$ws.Range("J51").Value = 43380
$ws.Range("K51").Value = 2090
$ws.Range("L51").Value = 43380
$ws.Range("M51").Value = -1354
$ws.Range("N51").Value = -44852

# Row 58
$ws.Range("H58").Value = 25641894
$ws.Range("I58").Value = 35714988
$ws.Range("J58").Value = 1285.7273
$ws.Range("K58").Value = 35714988
$ws.Range("L58").Value = 1285.7273
$ws.Range("M58").Value = -35714785
$ws.Range("N58").Value = -1691.7273

# Row 61
$ws.Range("H61").Value = 33057.5
$ws.Range("I61").Value = 2090
$ws.Range("J61").Value = 43380
$ws.Range("K61").Value = 2090
$ws.Range("L61").Value = 43380
$ws.Range("M61").Value = -1742
$ws.Range("N61").Value = -44076

# Row 122
$ws.Range("H122").Value = 19233270
$ws.Range("I122").Value = 27780546
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 83341638
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -83339188
$ws.Range("N122").Value = -10600

# Row 132
$ws.Range("H132").Value = 6290954
$ws.Range("I132").Value = 1319.8125
$ws.Range("J132").Value = 15875159
$ws.Range("K132").Value = 3959.4375
$ws.Range("L132").Value = 47625477
$ws.Range("M132").Value = -1429.4375
$ws.Range("N132").Value = -47630537

# Row 136
$ws.Range("H136").Value = 25641894
$ws.Range("I136").Value = 35714988
$ws.Range("J136").Value = 1285.7273
$ws.Range("K136").Value = 107144964
$ws.Range("L136").Value = 3857.1819
$ws.Range("M136").Value = -107142414
$ws.Range("N136").Value = -8957.1819

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 720.8
$ws.Range("I51").Value = 401
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 1203
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = -743
$ws.Range("N51").Value = -6920

# Row 107
$ws.Range("H107").Value = 2911.647
$ws.Range("I107").Value = 205.55556
$ws.Range("J107").Value = 3638.6567
$ws.Range("K107").Value = 616.66668
$ws.Range("L107").Value = 10915.9701
$ws.Range("M107").Value = 1303.33332
$ws.Range("N107").Value = -14755.9701

# Row 131
$ws.Range("H131").Value = 9473978
$ws.Range("I131").Value = 21742066
$ws.Range("J131").Value = 5132962
$ws.Range("K131").Value = 65226198
$ws.Range("L131").Value = 15398886
$ws.Range("M131").Value = -65221158
$ws.Range("N131").Value = -15408966

# Row 132
$ws.Range("H132").Value = 4741.2144
$ws.Range("I132").Value = 633.3333
$ws.Range("J132").Value = 6687.0527
$ws.Range("K132").Value = 5699.9997
$ws.Range("L132").Value = 60183.4743
$ws.Range("M132").Value = -3169.9997
$ws.Range("N132").Value = -65243.4743

# Row 140
$ws.Range("H140").Value = 7814107
$ws.Range("I140").Value = 10870663
$ws.Range("K140").Value = 32611989
$ws.Range("M140").Value = -32606809

# Row 141
$ws.Range("H141").Value = 2265.5
$ws.Range("I141").Value = 2009.6666
$ws.Range("K141").Value = 6028.9998
$ws.Range("M141").Value = -848.9997999999996

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 36015570
$ws.Range("I122").Value = 60023436
$ws.Range("J122").Value = 3765.1667
$ws.Range("K122").Value = 180070308
$ws.Range("L122").Value = 11295.5001
$ws.Range("M122").Value = -180067858
$ws.Range("N122").Value = -16195.5001

# Row 132
$ws.Range("H132").Value = 5385.0356
$ws.Range("I132").Value = 1484.762
$ws.Range("J132").Value = 17085.857
$ws.Range("K132").Value = 4454.286
$ws.Range("L132").Value = 51257.571
$ws.Range("M132").Value = -1924.286
$ws.Range("N132").Value = -56317.571

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1206.1538
$ws.Range("I68").Value = 1189.75
$ws.Range("K68").Value = 1189.75
$ws.Range("M68").Value = -440.75

# Row 71
$ws.Range("H71").Value = 1206.1538
$ws.Range("I71").Value = 1189.75
$ws.Range("K71").Value = 5948.75
$ws.Range("M71").Value = -2204.75

# Row 100
$ws.Range("H100").Value = 2618.1
$ws.Range("I100").Value = 1996.125
$ws.Range("J100").Value = 3032.75
$ws.Range("K100").Value = 1996.125
$ws.Range("L100").Value = 3032.75
$ws.Range("M100").Value = -1455.125
$ws.Range("N100").Value = -4114.75

# Row 101
$ws.Range("H101").Value = 19840.5
$ws.Range("J101").Value = 19840.5
$ws.Range("L101").Value = 19840.5
$ws.Range("N101").Value = -26330.5

# Row 103
$ws.Range("H103").Value = 29900
$ws.Range("J103").Value = 29900
$ws.Range("L103").Value = 29900
$ws.Range("N103").Value = -32244

# Row 132
$ws.Range("H132").Value = 22864224
$ws.Range("I132").Value = 42330500
$ws.Range("J132").Value = 12509.956
$ws.Range("K132").Value = 126991500
$ws.Range("L132").Value = 37529.868
$ws.Range("M132").Value = -126988970
$ws.Range("N132").Value = -42589.868

# Row 136
$ws.Range("H136").Value = 46587336
$ws.Range("I136").Value = 42331388
$ws.Range("J136").Value = 52635260
$ws.Range("K136").Value = 126994164
$ws.Range("L136").Value = 157905780
$ws.Range("M136").Value = -126991614
$ws.Range("N136").Value = -157910880

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 34475
$ws.Range("I132").Value = 51180.19
$ws.Range("J132").Value = 9417.214
$ws.Range("K132").Value = 153540.57
$ws.Range("L132").Value = 28251.642
$ws.Range("M132").Value = -151010.57
$ws.Range("N132").Value = -33311.642
